$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Venta - Plan"
